# Case_1_40 (380 kV) power-flow results: res_bus/vm_pu.xlsx
# The slack-bus voltage setpoint for this case was changed from 1.05 p.u. to
# 1.02 p.u.; re-running the load flow updates every bus voltage magnitude in
# rows 2-25 (24 time steps) for buses in columns B-F and I-N (column G is the
# external grid / reference bus and stays fixed at 1.0 p.u.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bus-voltage columns, left to right (column G is intentionally skipped).
$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

# New vm_pu results, one row per data row (sheet row 2 = time step 0, ...,
# sheet row 25 = time step 23), values aligned with $columns above.
$newVmPu = @{
    2 = @(1.02, 1.02008666070398, 1.034261023275327, 1.030171700409559, 1.038915499678274, 1.033509080966286, 1.02528575741646, 1.037061006112029, 1.032983483861769, 1.041702173788075, 1.012507370540426)
    3 = @(1.02, 1.021155830324076, 1.034790527288608, 1.031177595475883, 1.040148260606695, 1.033686467330915, 1.025991337788148, 1.03740040677212, 1.033797146900308, 1.042743930934674, 1.012740635244039)
    4 = @(1.02, 1.021845376943543, 1.035128666886632, 1.031825018975338, 1.040938105465185, 1.033795598235654, 1.026445102742728, 1.037614837095609, 1.034319587825441, 1.043409646010229, 1.012890625720032)
    5 = @(1.02, 1.022134720737394, 1.035269746425396, 1.032096371687901, 1.041268287298641, 1.033840123913095, 1.026635199619279, 1.037703741881792, 1.034538254932535, 1.043687513970806, 1.01295345578354)
    6 = @(1.02, 1.022183271175306, 1.035293371275693, 1.03214188477611, 1.041323616862526, 1.033847520616611, 1.026667078745182, 1.037718596550657, 1.034574913484265, 1.043734052169367, 1.012963992002775)
    7 = @(1.02, 1.021849245292398, 1.035130556220414, 1.031828648036589, 1.040942524705643, 1.033796198506435, 1.026447645439275, 1.037616029925324, 1.034322513461178, 1.043413366741228, 1.012891466144934)
    8 = @(1.02, 1.020448465455977, 1.034440899964539, 1.030512365656282, 1.03933374316892, 1.033570199253946, 1.025524791539767, 1.037176781438759, 1.033259306509276, 1.04205597734812, 1.012586400150322)
    9 = @(1.02, 1.017962491892162, 1.033191300838135, 1.02816624177109, 1.036438542379633, 1.033128687448099, 1.02387708395726, 1.03636306825862, 1.031354587844049, 1.039599690239053, 1.012041535903474)
    10 = @(1.02, 1.016293062297336, 1.032335138774957, 1.026583947588735, 1.034467349400126, 1.032805226359944, 1.022763938810997, 1.035793882773077, 1.030063529428932, 1.03791843290678, 1.01167331908334)
    11 = @(1.02, 1.015567242599783, 1.031958923136209, 1.0258944066485, 1.03360393688068, 1.032658244274816, 1.022278403564963, 1.035541072570493, 1.029499382103658, 1.037179945828323, 1.011512680791613)
    12 = @(1.02, 1.015297192276918, 1.031818353362765, 1.025637613188654, 1.033281731904073, 1.032602607056793, 1.022097518129406, 1.035446212162626, 1.029289058673161, 1.036904052001708, 1.011452831034523)
    13 = @(1.02, 1.015355139386082, 1.031848543485352, 1.025692726562464, 1.033350913819724, 1.032614588590354, 1.02213634300507, 1.035466603306138, 1.029334208891159, 1.036963304149364, 1.011465677251726)
    14 = @(1.02, 1.015544929342329, 1.031947320465446, 1.025873193672457, 1.0335773338877, 1.032653666542648, 1.022263462481764, 1.035533250880427, 1.029482012529699, 1.037157172793444, 1.011507737304775)
    15 = @(1.02, 1.015661805573681, 1.032008070650371, 1.025984296792834, 1.03371664028068, 1.032677605690945, 1.022341713810815, 1.035574187993385, 1.029572976476065, 1.037276411062534, 1.011533627787268)
    16 = @(1.02, 1.016341169991031, 1.032359991144896, 1.026629616949708, 1.034524442239873, 1.032814835116433, 1.022796087283963, 1.035810527058813, 1.030100861835409, 1.037967221926741, 1.011683954746405)
    17 = @(1.02, 1.0167665245403, 1.032579270332431, 1.027033227431947, 1.035028503662955, 1.032899061389043, 1.023080153767819, 1.035957075465264, 1.030430617882377, 1.038397733163747, 1.011777929068172)
    18 = @(1.02, 1.017014342861455, 1.032706642410045, 1.027268222810066, 1.035321562405718, 1.032947521476122, 1.023245504102349, 1.036041942196901, 1.030622466355763, 1.03864783155784, 1.011832627223095)
    19 = @(1.02, 1.017098794546182, 1.032749983178288, 1.027348278473906, 1.03542132677543, 1.032963931910208, 1.023301826621103, 1.036070775714601, 1.030687798346809, 1.038732937384243, 1.011851258339239)
    20 = @(1.02, 1.016720917419459, 1.032555798556416, 1.026989967728409, 1.034974521171172, 1.032890093779005, 1.023049711411041, 1.035941415569771, 1.030395289180356, 1.038351648084304, 1.011767858468235)
    21 = @(1.02, 1.015489053362996, 1.031918255936345, 1.025820069081123, 1.033510700192068, 1.032642187825908, 1.022226043800826, 1.035513651219849, 1.029438509459661, 1.037100127202935, 1.011495356691895)
    22 = @(1.02, 1.014711933073946, 1.031512624571788, 1.025080643299531, 1.03258168117708, 1.03248029251175, 1.021705067446951, 1.035239169909103, 1.028832462648646, 1.036304059627949, 1.011322972823664)
    23 = @(1.02, 1.015124147625916, 1.03172811137241, 1.025472995475404, 1.03307499640879, 1.032566688209097, 1.021981542737744, 1.035385202318755, 1.029154166479358, 1.036726944619314, 1.01141445695669)
    24 = @(1.02, 1.016741526191386, 1.032566406078993, 1.02700951623605, 1.034998916476587, 1.032894147919266, 1.023063468055032, 1.035948493495287, 1.030411254222527, 1.038372475073713, 1.011772409295737)
    25 = @(1.02, 1.018607290259752, 1.033518421075894, 1.028775958886752, 1.037194218860303, 1.033247958362286, 1.024305626268873, 1.036578138019562, 1.031850727598841, 1.04024234133276, 1.012183267885265)
}

foreach ($row in $newVmPu.Keys) {
    $rowValues = $newVmPu[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $rowValues[$i]
    }
}

